$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 176, 177, 178 hold three Azerbaijan Premier League matches scheduled
# for the same date/round. The data (everything except the id column A,
# which just reflects its own row position) needs to be cyclically rotated:
#   new row176 = old row177 data
#   new row177 = old row178 data
#   new row178 = old row176 data

# Capture the existing values (columns B..AD) before any writes.
$row176 = $ws.Range("B176:AD176").Value2
$row177 = $ws.Range("B177:AD177").Value2
$row178 = $ws.Range("B178:AD178").Value2

$ws.Range("B176:AD176").Value2 = $row177
$ws.Range("B177:AD177").Value2 = $row178
$ws.Range("B178:AD178").Value2 = $row176
